$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Vega Monumental Concepción - Acelga".
# It belongs chronologically right after the existing row 50 (date 2021-06-18,
# serial 44365) and before the former row 51 (date 2021-07-28, serial 44405),
# so insert a fresh row at position 51 - this pushes every following row down
# by one (old row 51 -> new row 52, ..., old row 170 -> new row 171) exactly
# as the diff shows.
$ws.Rows("51").Insert()

# Populate the newly inserted row 51 with the new record's data.
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 44497
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112009
$ws.Range("G51").Value = "Acelga"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 250
$ws.Range("K51").Value = 600
$ws.Range("L51").Value = 650
$ws.Range("M51").Value = 630
$ws.Range("N51").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 420
$ws.Range("Q51").Value = 1.5
$ws.Range("R51").Value = "Hortaliza"
